$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric data in B2:E13 to the nearest integer, in place.
$rng = $ws.Range("B2:E13")
foreach ($cell in $rng.Cells) {
    $val = $cell.Value()
    if ($null -ne $val) {
        $cell.Value = $excel.WorksheetFunction.Round($val, 0)
    }
}
